# Revlon UK Email triggering Code Latest changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order below matches the order new strings were appended to sharedStrings.xml
# (index 83..94) in the authored workbook - mirror that exact entry order here.

# 83 HeaderNames
$ws.Range("AD1").Value = "HeaderNames"

# 84 HeaderLinks
$ws.Range("A18").Value = "HeaderLinks"

# 85 One-Step,Straighteners,Dryers,Hair Stylers
$ws.Range("AD18").Value = "One-Step,Straighteners,Dryers,Hair Stylers"

# 86 ForgotPassword
$ws.Range("A19").Value = "ForgotPassword"

# (reuse of existing shared strings 22 / 21)
$ws.Range("C19").Value = "Harish!123"
$ws.Range("F19").Value = "harish.chiruvella1@gmail.com"
$ws.Hyperlinks.Add($ws.Range("F19"), "mailto:harish.chiruvella1@gmail.com") | Out-Null
$ws.Range("F19").Style = "Hyperlink"

$ws.Range("Y19").Font.Bold = $false
$ws.Range("Z19").Font.Bold = $false
$ws.Range("AA19").Font.Bold = $false

# 87 Revlon Hair Tools <Revlon@r1.dotdigital-email.com>
$ws.Range("AE19").Value = "Revlon Hair Tools <Revlon@r1.dotdigital-email.com>"
# 88 Harish Chiruvella <harish.chiruvella1@gmail.com>
$ws.Range("AF19").Value = "Harish Chiruvella <harish.chiruvella1@gmail.com>"
# 89 Your Password Reset Request
$ws.Range("AG19").Value = "Your Password Reset Request"
# 90 Set a New Password
$ws.Range("AH19").Value = "Set a New Password"

# 91 fromAddress
$ws.Range("AE1").Value = "fromAddress"
# 92 toAddress
$ws.Range("AF1").Value = "toAddress"
# 93 subject
$ws.Range("AG1").Value = "subject"
# 94 content
$ws.Range("AH1").Value = "content"

$ws.Range("AD1:AH1").Interior.Color = $ws.Range("AC1").Interior.Color

# Auto-fit the newly populated columns (AD..AH) to their content, mirroring the
# bestFit column widths Excel records after typing this data.
for ($col = 30; $col -le 34; $col++) {
    $ws.Columns.Item($col).AutoFit() | Out-Null
}

# Reselect range to match final workbook selection state
$ws.Range("AE1:AH1").Select()

$wb.Save()
